# Insert a new multiple-choice question as row 15, pushing the existing
# rows 15-40 down to rows 16-41 (new dimension becomes A1:K41).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 15 (shifts rows 15..40 down to 16..41)
$ws.Rows("15:15").Insert()

# Populate the newly inserted row 15 with the new question
$ws.Range("A15").Value = "mc"
$ws.Range("B15").Value = "A 28-year-old female figure skater presents several weeks after having sustained an injury to her left breast. She has a painful mass shown in the upper outer quadrant. Skin retraction is noticed, and a hard mass, 3–4 cm in diameter, can easily be palpated. What is the most likely diagnosis?"
$ws.Range("C15").Value = "Infiltrating carcinoma."
$ws.Range("D15").Value = "Breast abscess."
$ws.Range("E15").Value = "Hematoma."
$ws.Range("F15").Value = "Fat necrosis."
$ws.Range("G15").Value = "Sclerosing adenosis."
$ws.Range("H15").Value = 4
